$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 79; this shifts rows 79-126 down to 80-127
# (and the sheet dimension grows from A1:R126 to A1:R127 automatically).
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly record.
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 44729
$ws.Range("D79").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112052
$ws.Range("G79").Value = "Albahaca"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 150
$ws.Range("K79").Value = 6000
$ws.Range("L79").Value = 7000
$ws.Range("M79").Value = 6333
$ws.Range("N79").Value = "$/paquete"
$ws.Range("O79").Value = "Región de Arica y Parinacota"
$ws.Range("P79").Value = 6333
$ws.Range("Q79").Value = 1
$ws.Range("R79").Value = "Hortaliza"
